# Apply odds updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 5
$ws.Range("I2").Value = 1.7
$ws.Range("Q2").Value = 2.05
$ws.Range("R2").Value = 1.85
$ws.Range("AC2").Value = 9.5
$ws.Range("AD2").Value = 7
$ws.Range("AK2").Value = 13
$ws.Range("AX2").Value = 9

# Row 4
$ws.Range("G4").Value = 3.4
$ws.Range("I4").Value = 2.25
$ws.Range("L4").Value = 2.88
$ws.Range("W4").Value = 10
$ws.Range("X4").Value = 17
$ws.Range("AC4").Value = 9
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 13
$ws.Range("AO4").Value = 19

# Row 6
$ws.Range("O6").Value = 1.14
$ws.Range("P6").Value = 5
$ws.Range("Q6").Value = 1.5
$ws.Range("R6").Value = 2.5

# Row 7
$ws.Range("U7").Value = 1.67

# Row 8
$ws.Range("U8").Value = 1.53
$ws.Range("V8").Value = 2.38
